# Auto commit at 2025-12-15 11:51:06.65
# Update Metrics values and let dependent formulas on "today" sheet recalc.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update source values in column B ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 190007.48
$metrics.Range("B3").Value = 162819.75000000003
$metrics.Range("B4").Value = 58277.289999999994
$metrics.Range("B5").Value = 7736
$metrics.Range("B6").Value = 5392714.5900000008
$metrics.Range("B7").Value = 4563172.7100000009
$metrics.Range("B8").Value = 1590234.1700000004
$metrics.Range("B9").Value = 210443
$metrics.Range("B10").Value = 33858095.579999991
$metrics.Range("B11").Value = 31838447.870000001
$metrics.Range("B12").Value = 11871956.209999995
$metrics.Range("B13").Value = 1308073

# Restore the selection recorded on the Metrics sheet view
$metrics.Activate()
$metrics.Range("F20:F21").Select()

# --- today sheet: selection moved, formulas/TODAY() recalc automatically ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("D5").Select()

$excel.Calculate()
